# Update "想去人数" (interest count) figures in column F for the
# "展览" and "全部类型" worksheets, reflecting newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 5539
$ws1.Range("F5").Value  = 59
$ws1.Range("F7").Value  = 659
$ws1.Range("F13").Value = 5085
$ws1.Range("F15").Value = 238
$ws1.Range("F16").Value = 213
$ws1.Range("F17").Value = 26
$ws1.Range("F18").Value = 9
$ws1.Range("F20").Value = 4369
$ws1.Range("F22").Value = 1157
$ws1.Range("F24").Value = 64
$ws1.Range("F27").Value = 171
$ws1.Range("F32").Value = 12
$ws1.Range("F33").Value = 41

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 5539
$ws4.Range("F6").Value  = 59
$ws4.Range("F8").Value  = 659
$ws4.Range("F14").Value = 5085
$ws4.Range("F16").Value = 238
$ws4.Range("F17").Value = 213
$ws4.Range("F18").Value = 26
$ws4.Range("F19").Value = 9
$ws4.Range("F21").Value = 4369
$ws4.Range("F23").Value = 1157
$ws4.Range("F25").Value = 64
$ws4.Range("F28").Value = 171
$ws4.Range("F33").Value = 12
$ws4.Range("F34").Value = 41

$wb.Save()
